$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows appended to the bottom of the table (dates as Excel serials)
$rows = @(
    @{ Row = 245; A = 44319; B = 0; C = 1; D = 40.79967360261118 },
    @{ Row = 246; A = 44320; B = 1; C = 2; D = 81.59934720522236 },
    @{ Row = 247; A = 44321; B = 0; C = 2; D = 81.59934720522236 }
)

$lastRow = 244

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A keeps the same style (border/bold/date numfmt) as the row above it
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowNum, 1).Value = $r.A

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    $lastRow = $rowNum
}
